$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force D-column target cells to remain text (they hold numeric-looking strings)
$dCells = @("D2","D3","D4","D5","D7","D8","D9","D10","D11","D13","D14","D15","D16","D17","D18","D20","D21","D23","D24","D25","D26","D27","D28","D29","D31","D32","D33","D34","D35","D36","D37","D39","D40","D41","D42","D43","D44","D45","D46","D47","D49","D50","D51")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "30.043.75"
$ws.Range("E2").Value = "  +3.93%  "
$ws.Range("D3").Value = "1.895.09"
$ws.Range("E3").Value = "  +4.21%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").Value = "248.40"
$ws.Range("E5").Value = "  +0.61%  "
$ws.Range("E6").Value = "  +0.19%  "
$ws.Range("D7").Value = "0.4959"
$ws.Range("E7").Value = "  +0.39%  "
$ws.Range("D8").Value = "44.91"
$ws.Range("E8").Value = "  +2.77%  "
$ws.Range("D9").Value = "0.2946"
$ws.Range("E9").Value = "  +5.94%  "
$ws.Range("D10").Value = "0.06629"
$ws.Range("E10").Value = "  +3.60%  "
$ws.Range("D11").Value = "1.899.10"
$ws.Range("E11").Value = "  +4.51%  "
$ws.Range("E12").Value = "  +1.66%  "
$ws.Range("D13").Value = "0.07230"
$ws.Range("E13").Value = "  +2.34%  "
$ws.Range("D14").Value = "0.6768"
$ws.Range("E14").Value = "  +5.26%  "
$ws.Range("D15").Value = "86.16"
$ws.Range("E15").Value = "  +2.40%  "
$ws.Range("D16").Value = "4.856"
$ws.Range("E16").Value = "  +3.98%  "
$ws.Range("D17").Value = "30.044.45"
$ws.Range("E17").Value = "  +3.88%  "
$ws.Range("D18").Value = "0.000007916"
$ws.Range("E18").Value = "  +8.20%  "
$ws.Range("E19").Value = "  +0.16%  "
$ws.Range("D20").Value = "12.89"
$ws.Range("E20").Value = "  +5.54%  "
$ws.Range("D21").Value = "2.145.36"
$ws.Range("E21").Value = "  +5.11%  "
$ws.Range("E22").Value = "  +0.26%  "
$ws.Range("D23").Value = "4.769"
$ws.Range("E23").Value = "  +4.61%  "
$ws.Range("D24").Value = "5.662"
$ws.Range("E24").Value = "  +5.79%  "
$ws.Range("D25").Value = "9.212"
$ws.Range("E25").Value = "  +4.30%  "
$ws.Range("D26").Value = "147.43"
$ws.Range("E26").Value = "  +1.85%  "
$ws.Range("D27").Value = "131.63"
$ws.Range("E27").Value = "  +2.29%  "
$ws.Range("D28").Value = "16.78"
$ws.Range("E28").Value = "  +2.30%  "
$ws.Range("D29").Value = "1.963"
$ws.Range("E29").Value = "  +4.38%  "
$ws.Range("E30").Value = "  -1.46%  "
$ws.Range("D31").Value = "4.220"
$ws.Range("E31").Value = "  +2.12%  "
$ws.Range("D32").Value = "0.08742"
$ws.Range("E32").Value = "  +4.49%  "
$ws.Range("D33").Value = "3.938"
$ws.Range("E33").Value = "  +3.96%  "
$ws.Range("D34").Value = "0.05092"
$ws.Range("E34").Value = "  +3.51%  "
$ws.Range("D35").Value = "1.123"
$ws.Range("E35").Value = "  +2.43%  "
$ws.Range("D36").Value = "0.7012"
$ws.Range("E36").Value = "  +3.94%  "
$ws.Range("D37").Value = "2.681"
$ws.Range("E37").Value = "  -0.54%  "
$ws.Range("D39").Value = "2.224"
$ws.Range("E39").Value = "  -3.31%  "
$ws.Range("D40").Value = "0.9527"
$ws.Range("E40").Value = "  +0.49%  "
$ws.Range("D41").Value = "0.01661"
$ws.Range("E41").Value = "  +4.79%  "
$ws.Range("D42").Value = "5.979"
$ws.Range("E42").Value = "  -2.76%  "

# Row 43 and 44: swap TheSandbox and PaxDollar entries
$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D43").Value = "1.000"
$ws.Range("E43").Value = "  +0.07%  "
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").Value = "0.4222"
$ws.Range("E44").Value = "  +3.42%  "

$ws.Range("D45").Value = "103.07"
$ws.Range("E45").Value = "  +2.56%  "
$ws.Range("D46").Value = "7.469"
$ws.Range("E46").Value = "  +4.06%  "
$ws.Range("D47").Value = "0.1258"
$ws.Range("E47").Value = "  +2.89%  "
$ws.Range("E48").Value = "  +4.27%  "
$ws.Range("D49").Value = "32.80"
$ws.Range("E49").Value = "  +3.63%  "
$ws.Range("D50").Value = "8.282"
$ws.Range("E50").Value = "  +2.35%  "
$ws.Range("D51").Value = "0.3734"
$ws.Range("E51").Value = "  +3.43%  "

# Restore default style on the D-column cells we touched (removes temp text format)
foreach ($addr in $dCells) {
    $ws.Range($addr).Style = "Normal"
}
